$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 216.72223
$ws.Range("I33").Value = 109.916664
$ws.Range("K33").Value = 109.916664
$ws.Range("M33").Value = 119.083336

$ws.Range("H98").Value = 1119281.6
$ws.Range("I98").Value = 1119281.6
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1119281.6
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1117783.6
$ws.Range("N98").ClearContents()

$ws.Range("H112").Value = 16043820
$ws.Range("J112").Value = 16043820
$ws.Range("L112").Value = 48131460
$ws.Range("N112").Value = -48133676

$ws.Range("H122").Value = 1119281.6
$ws.Range("I122").Value = 1119281.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3357844.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3355394.8
$ws.Range("N122").ClearContents()

$ws.Range("H129").Value = 921.62
$ws.Range("I129").Value = 380.1111
$ws.Range("J129").Value = 1040.4878
$ws.Range("K129").Value = 1140.3333
$ws.Range("L129").Value = 3121.463400000001
$ws.Range("M129").Value = 3859.6667
$ws.Range("N129").Value = -13121.4634

$ws.Range("H138").Value = 5664302.5
$ws.Range("I138").Value = 1536828.9
$ws.Range("J138").Value = 6851658
$ws.Range("K138").Value = 4610486.699999999
$ws.Range("L138").Value = 20554974
$ws.Range("M138").Value = -4605346.699999999
$ws.Range("N138").Value = -20565254

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2428.0984
$ws.Range("I32").Value = 1625.7091
$ws.Range("K32").Value = 1625.7091
$ws.Range("M32").Value = -1338.7091

$ws.Range("H74").Value = 6729.095
$ws.Range("I74").Value = 1105.8422
$ws.Range("J74").Value = 60150
$ws.Range("K74").Value = 1105.8422
$ws.Range("L74").Value = 60150
$ws.Range("M74").Value = -231.8422
$ws.Range("N74").Value = -61898

$ws.Range("H77").Value = 6729.095
$ws.Range("I77").Value = 1105.8422
$ws.Range("J77").Value = 60150
$ws.Range("K77").Value = 5529.211
$ws.Range("L77").Value = 300750
$ws.Range("M77").Value = -1161.211
$ws.Range("N77").Value = -309486

$ws.Range("H123").Value = 32950.668
$ws.Range("J123").Value = 32950.668
$ws.Range("L123").Value = 32950.668
$ws.Range("N123").Value = -42750.668

$ws.Range("H139").Value = 47079.6
$ws.Range("J139").Value = 47079.6
$ws.Range("L139").Value = 47079.6
$ws.Range("N139").Value = -57359.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 3272.7273
$ws.Range("J15").Value = 3272.7273
$ws.Range("L15").Value = 3272.7273
$ws.Range("N15").Value = -3726.7273

$ws.Range("H99").Value = 1692.5
$ws.Range("I99").Value = 2148
$ws.Range("J99").Value = 933.3333
$ws.Range("K99").Value = 2148
$ws.Range("L99").Value = 933.3333
$ws.Range("M99").Value = -650
$ws.Range("N99").Value = -3929.3333

$ws.Range("H133").Value = 54000
$ws.Range("J133").Value = 54000
$ws.Range("L133").Value = 54000
$ws.Range("N133").Value = -64120

$ws.Range("H134").Value = 3363.4138
$ws.Range("I134").Value = 1940.45
$ws.Range("K134").Value = 5821.35
$ws.Range("M134").Value = -3286.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 274
$ws.Range("I33").Value = 249.71428
$ws.Range("J33").Value = 330.66666
$ws.Range("K33").Value = 1498.28568
$ws.Range("L33").Value = 1983.99996
$ws.Range("M33").Value = -1215.28568
$ws.Range("N33").Value = -2549.99996

$ws.Range("H116").Value = 1290
$ws.Range("I116").Value = 886.6667
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2660.0001
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = 781.9998999999998
$ws.Range("N116").Value = -14384

$ws.Range("H131").Value = 3322.0894
$ws.Range("J131").Value = 3492.2075
$ws.Range("L131").Value = 10476.6225
$ws.Range("N131").Value = -20556.6225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 42000
$ws.Range("J137").Value = 42000
$ws.Range("L137").Value = 42000
$ws.Range("N137").Value = -52200

$ws.Range("H138").Value = 48666.668
$ws.Range("J138").Value = 48666.668
$ws.Range("L138").Value = 48666.668
$ws.Range("N138").Value = -58946.668

$ws.Range("H139").Value = 54980
$ws.Range("J139").Value = 54980
$ws.Range("L139").Value = 54980
$ws.Range("N139").Value = -65260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 19000
$ws.Range("I4").Value = 19000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -18887
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 19000
$ws.Range("I28").Value = 19000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 19000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -18768
$ws.Range("N28").ClearContents()

$ws.Range("H37").Value = 19000
$ws.Range("I37").Value = 19000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 19000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -18893
$ws.Range("N37").ClearContents()

$ws.Range("H68").Value = 1380.5
$ws.Range("I68").Value = 1422.7778
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1422.7778
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -673.7778000000001
$ws.Range("N68").Value = -2498

$ws.Range("H71").Value = 1380.5
$ws.Range("I71").Value = 1422.7778
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 7113.889
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -3369.889
$ws.Range("N71").Value = -12488

$ws.Range("H122").Value = 3494.1875
$ws.Range("I122").Value = 2636
$ws.Range("J122").Value = 3692.2307
$ws.Range("K122").Value = 7908
$ws.Range("L122").Value = 11076.6921
$ws.Range("M122").Value = -5458
$ws.Range("N122").Value = -15976.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 18575508
$ws.Range("I136").Value = 22290136
$ws.Range("J136").Value = 2365
$ws.Range("K136").Value = 66870408
$ws.Range("L136").Value = 7095
$ws.Range("M136").Value = -66867858
$ws.Range("N136").Value = -12195
